$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("C13").Value = "Dec 2023 (24/06/24)"
$ws.Range("C13").Select()
